$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Cells.Item(29, 8).Value = 100001810   # H29: 83334856 -> 100001810
$ws.Cells.Item(29, 10).Value = 7999   # J29: 4042 -> 7999
$ws.Cells.Item(29, 12).Value = 23997   # L29: 12126 -> 23997
$ws.Cells.Item(29, 14).Value = -24559   # N29: -12688 -> -24559
# Row 74
$ws.Cells.Item(74, 8).Value = 15290.1   # H74: 14345.091 -> 15290.1
$ws.Cells.Item(74, 9).Value = 15290.1   # I74: 14345.091 -> 15290.1
$ws.Cells.Item(74, 11).Value = 15290.1   # K74: 14345.091 -> 15290.1
$ws.Cells.Item(74, 13).Value = -14354.1   # M74: -13409.091 -> -14354.1
# Row 77
$ws.Cells.Item(77, 8).Value = 15290.1   # H77: 14345.091 -> 15290.1
$ws.Cells.Item(77, 9).Value = 15290.1   # I77: 14345.091 -> 15290.1
$ws.Cells.Item(77, 11).Value = 76450.5   # K77: 71725.455 -> 76450.5
$ws.Cells.Item(77, 13).Value = -71770.5   # M77: -67045.455 -> -71770.5
# Row 95
$ws.Cells.Item(95, 8).Value = 45687.168   # H95: 49333 -> 45687.168
$ws.Cells.Item(95, 10).Value = 45687.168   # J95: 49333 -> 45687.168
$ws.Cells.Item(95, 12).Value = 45687.168   # L95: 49333 -> 45687.168
$ws.Cells.Item(95, 14).Value = -51179.168   # N95: -54825 -> -51179.168
# Row 98
$ws.Cells.Item(98, 8).Value = 7367.8335   # H98: 8001.864 -> 7367.8335
$ws.Cells.Item(98, 9).Value = 10013.412   # I98: 11296.066 -> 10013.412
$ws.Cells.Item(98, 11).Value = 10013.412   # K98: 11296.066 -> 10013.412
$ws.Cells.Item(98, 13).Value = -8515.412   # M98: -9798.066000000001 -> -8515.412
# Row 112
$ws.Cells.Item(112, 8).Value = 2568.5   # H112: 2416.4443 -> 2568.5
$ws.Cells.Item(112, 9).Value = 1341.5   # I112: 1321.2858 -> 1341.5
$ws.Cells.Item(112, 11).Value = 4024.5   # K112: 3963.8574 -> 4024.5
$ws.Cells.Item(112, 13).Value = -2916.5   # M112: -2855.8574 -> -2916.5
# Row 122
$ws.Cells.Item(122, 8).Value = 7367.8335   # H122: 8001.864 -> 7367.8335
$ws.Cells.Item(122, 9).Value = 10013.412   # I122: 11296.066 -> 10013.412
$ws.Cells.Item(122, 11).Value = 30040.236   # K122: 33888.198 -> 30040.236
$ws.Cells.Item(122, 13).Value = -27590.236   # M122: -31438.198 -> -27590.236
# Row 132
$ws.Cells.Item(132, 8).Value = 1756.1177   # H132: 1777.25 -> 1756.1177
$ws.Cells.Item(132, 9).Value = 1436.1072   # I132: 1458.2963 -> 1436.1072
$ws.Cells.Item(132, 10).Value = 3249.5   # J132: 3499.6 -> 3249.5
$ws.Cells.Item(132, 11).Value = 4308.321599999999   # K132: 4374.8889 -> 4308.321599999999
$ws.Cells.Item(132, 12).Value = 9748.5   # L132: 10498.8 -> 9748.5
$ws.Cells.Item(132, 13).Value = -1778.321599999999   # M132: -1844.8889 -> -1778.321599999999
$ws.Cells.Item(132, 14).Value = -14808.5   # N132: -15558.8 -> -14808.5
# Row 137
$ws.Cells.Item(137, 8).Value = 1405.1333   # H137: 1437.3448 -> 1405.1333
$ws.Cells.Item(137, 9).Value = 1475.125   # I137: 1445.68 -> 1475.125
$ws.Cells.Item(137, 10).Value = 1125.1666   # J137: 1385.25 -> 1125.1666
$ws.Cells.Item(137, 11).Value = 4425.375   # K137: 4337.04 -> 4425.375
$ws.Cells.Item(137, 12).Value = 3375.4998   # L137: 4155.75 -> 3375.4998
$ws.Cells.Item(137, 13).Value = -1875.375   # M137: -1787.04 -> -1875.375
$ws.Cells.Item(137, 14).Value = -8475.4998   # N137: -9255.75 -> -8475.4998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 2304.2344   # H32: 2249.3635 -> 2304.2344
$ws.Cells.Item(32, 9).Value = 2304.2344   # I32: 2249.3635 -> 2304.2344
$ws.Cells.Item(32, 11).Value = 2304.2344   # K32: 2249.3635 -> 2304.2344
$ws.Cells.Item(32, 13).Value = -2017.2344   # M32: -1962.3635 -> -2017.2344
# Row 45
$ws.Cells.Item(45, 8).Value = 1083.3334   # H45: 1125 -> 1083.3334
$ws.Cells.Item(45, 10).Value = 1000   # J45: 0 -> 1000
$ws.Cells.Item(45, 12).Value = 1000   # L45: 0 -> 1000
$ws.Cells.Item(45, 14).Value = -1754   # N45: None -> -1754
# Row 63
$ws.Cells.Item(63, 8).Value = 5752.3   # H63: 6719.4 -> 5752.3
$ws.Cells.Item(63, 9).Value = 3154.7778   # I63: 3424.75 -> 3154.7778
$ws.Cells.Item(63, 10).Value = 7877.5454   # J63: 8269.823 -> 7877.5454
$ws.Cells.Item(63, 11).Value = 3154.7778   # K63: 3424.75 -> 3154.7778
$ws.Cells.Item(63, 12).Value = 7877.5454   # L63: 8269.823 -> 7877.5454
$ws.Cells.Item(63, 13).Value = -2468.7778   # M63: -2738.75 -> -2468.7778
$ws.Cells.Item(63, 14).Value = -9249.545399999999   # N63: -9641.823 -> -9249.545399999999
# Row 66
$ws.Cells.Item(66, 8).Value = 5752.3   # H66: 6719.4 -> 5752.3
$ws.Cells.Item(66, 9).Value = 3154.7778   # I66: 3424.75 -> 3154.7778
$ws.Cells.Item(66, 10).Value = 7877.5454   # J66: 8269.823 -> 7877.5454
$ws.Cells.Item(66, 11).Value = 15773.889   # K66: 17123.75 -> 15773.889
$ws.Cells.Item(66, 12).Value = 39387.727   # L66: 41349.11500000001 -> 39387.727
$ws.Cells.Item(66, 13).Value = -12341.889   # M66: -13691.75 -> -12341.889
$ws.Cells.Item(66, 14).Value = -46251.727   # N66: -48213.11500000001 -> -46251.727
# Row 122
$ws.Cells.Item(122, 8).Value = 1510.1538   # H122: 1562.12 -> 1510.1538
$ws.Cells.Item(122, 9).Value = 1173.2858   # I122: 1221.4 -> 1173.2858
$ws.Cells.Item(122, 11).Value = 3519.8574   # K122: 3664.2 -> 3519.8574
$ws.Cells.Item(122, 13).Value = -1069.8574   # M122: -1214.2 -> -1069.8574

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Cells.Item(82, 8).Value = 19849.143   # H82: 27741 -> 19849.143
$ws.Cells.Item(82, 9).Value = 13157.333   # I82: 14937.4 -> 13157.333
$ws.Cells.Item(82, 10).Value = 60000   # J82: 59750 -> 60000
$ws.Cells.Item(82, 11).Value = 13157.333   # K82: 14937.4 -> 13157.333
$ws.Cells.Item(82, 12).Value = 60000   # L82: 59750 -> 60000
$ws.Cells.Item(82, 13).Value = -12774.333   # M82: -14554.4 -> -12774.333
$ws.Cells.Item(82, 14).Value = -60766   # N82: -60516 -> -60766
# Row 85
$ws.Cells.Item(85, 8).Value = 19849.143   # H85: 27741 -> 19849.143
$ws.Cells.Item(85, 9).Value = 13157.333   # I85: 14937.4 -> 13157.333
$ws.Cells.Item(85, 10).Value = 60000   # J85: 59750 -> 60000
$ws.Cells.Item(85, 11).Value = 13157.333   # K85: 14937.4 -> 13157.333
$ws.Cells.Item(85, 12).Value = 60000   # L85: 59750 -> 60000
$ws.Cells.Item(85, 13).Value = -11831.333   # M85: -13611.4 -> -11831.333
$ws.Cells.Item(85, 14).Value = -62652   # N85: -62402 -> -62652
# Row 86
$ws.Cells.Item(86, 8).Value = 33337180   # H86: 35718336 -> 33337180
$ws.Cells.Item(86, 9).Value = 3640.5833   # I86: 3880.6365 -> 3640.5833
$ws.Cells.Item(86, 11).Value = 3640.5833   # K86: 3880.6365 -> 3640.5833
$ws.Cells.Item(86, 13).Value = -2517.5833   # M86: -2757.6365 -> -2517.5833
# Row 89
$ws.Cells.Item(89, 8).Value = 33337180   # H89: 35718336 -> 33337180
$ws.Cells.Item(89, 9).Value = 3640.5833   # I89: 3880.6365 -> 3640.5833
$ws.Cells.Item(89, 11).Value = 18202.9165   # K89: 19403.1825 -> 18202.9165
$ws.Cells.Item(89, 13).Value = -12586.9165   # M89: -13787.1825 -> -12586.9165

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 4933   # H31: 6048.231 -> 4933
$ws.Cells.Item(31, 9).Value = 1690.4286   # I31: 1876.1428 -> 1690.4286
$ws.Cells.Item(31, 10).Value = 12499   # J31: 10915.667 -> 12499
$ws.Cells.Item(31, 11).Value = 1690.4286   # K31: 1876.1428 -> 1690.4286
$ws.Cells.Item(31, 12).Value = 12499   # L31: 10915.667 -> 12499
$ws.Cells.Item(31, 13).Value = -1395.4286   # M31: -1581.1428 -> -1395.4286
$ws.Cells.Item(31, 14).Value = -13089   # N31: -11505.667 -> -13089
# Row 34
$ws.Cells.Item(34, 8).Value = 4933   # H34: 6048.231 -> 4933
$ws.Cells.Item(34, 9).Value = 1690.4286   # I34: 1876.1428 -> 1690.4286
$ws.Cells.Item(34, 10).Value = 12499   # J34: 10915.667 -> 12499
$ws.Cells.Item(34, 11).Value = 1690.4286   # K34: 1876.1428 -> 1690.4286
$ws.Cells.Item(34, 12).Value = 12499   # L34: 10915.667 -> 12499
$ws.Cells.Item(34, 13).Value = -1488.4286   # M34: -1674.1428 -> -1488.4286
$ws.Cells.Item(34, 14).Value = -12903   # N34: -11319.667 -> -12903
# Row 59
$ws.Cells.Item(59, 8).Value = 38644.11   # H59: 42514 -> 38644.11
$ws.Cells.Item(59, 9).Value = 34250   # I59: 35000 -> 34250
$ws.Cells.Item(59, 10).Value = 39899.57   # J59: 43766.332 -> 39899.57
$ws.Cells.Item(59, 11).Value = 34250   # K59: 35000 -> 34250
$ws.Cells.Item(59, 12).Value = 39899.57   # L59: 43766.332 -> 39899.57
$ws.Cells.Item(59, 13).Value = -33105   # M59: -33855 -> -33105
$ws.Cells.Item(59, 14).Value = -42189.57   # N59: -46056.332 -> -42189.57
# Row 97
$ws.Cells.Item(97, 8).Value = 78978   # H97: 90796 -> 78978
$ws.Cells.Item(97, 9).Value = 41000   # I97: 67000 -> 41000
$ws.Cells.Item(97, 10).Value = 104296.664   # J97: 96745 -> 104296.664
$ws.Cells.Item(97, 11).Value = 41000   # K97: 67000 -> 41000
$ws.Cells.Item(97, 12).Value = 104296.664   # L97: 96745 -> 104296.664
$ws.Cells.Item(97, 13).Value = -40009   # M97: -66009 -> -40009
$ws.Cells.Item(97, 14).Value = -106278.664   # N97: -98727 -> -106278.664
# Row 99
$ws.Cells.Item(99, 8).Value = 1254593.8   # H99: 1254731.1 -> 1254593.8
$ws.Cells.Item(99, 9).Value = 1433550   # I99: 1433707 -> 1433550
$ws.Cells.Item(99, 11).Value = 1433550   # K99: 1433707 -> 1433550
$ws.Cells.Item(99, 13).Value = -1432052   # M99: -1432209 -> -1432052
# Row 104
$ws.Cells.Item(104, 8).Value = 59997.5   # H104: 0 -> 59997.5
$ws.Cells.Item(104, 10).Value = 59997.5   # J104: 0 -> 59997.5
$ws.Cells.Item(104, 12).Value = 59997.5   # L104: 0 -> 59997.5
$ws.Cells.Item(104, 14).Value = -65239.5   # N104: None -> -65239.5
# Row 126
$ws.Cells.Item(126, 8).Value = 1254593.8   # H126: 1254731.1 -> 1254593.8
$ws.Cells.Item(126, 9).Value = 1433550   # I126: 1433707 -> 1433550
$ws.Cells.Item(126, 11).Value = 4300650   # K126: 4301121 -> 4300650
$ws.Cells.Item(126, 13).Value = -4298180   # M126: -4298651 -> -4298180

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Cells.Item(40, 8).Value = 65.833336   # H40: 68.75 -> 65.833336
$ws.Cells.Item(40, 9).Value = 48.333332   # I40: 50 -> 48.333332
$ws.Cells.Item(40, 10).Value = 83.333336   # J40: 87.5 -> 83.333336
$ws.Cells.Item(40, 11).Value = 193.333328   # K40: 200 -> 193.333328
$ws.Cells.Item(40, 12).Value = 333.333344   # L40: 350 -> 333.333344
$ws.Cells.Item(40, 13).Value = -124.333328   # M40: -131 -> -124.333328
$ws.Cells.Item(40, 14).Value = -471.333344   # N40: -488 -> -471.333344

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 4913.8184   # H80: 5076.381 -> 4913.8184
$ws.Cells.Item(80, 10).Value = 4883.4165   # J80: 5191 -> 4883.4165
$ws.Cells.Item(80, 12).Value = 4883.4165   # L80: 5191 -> 4883.4165
$ws.Cells.Item(80, 14).Value = -6879.4165   # N80: -7187 -> -6879.4165
# Row 83
$ws.Cells.Item(83, 8).Value = 4913.8184   # H83: 5076.381 -> 4913.8184
$ws.Cells.Item(83, 10).Value = 4883.4165   # J83: 5191 -> 4883.4165
$ws.Cells.Item(83, 12).Value = 24417.0825   # L83: 25955 -> 24417.0825
$ws.Cells.Item(83, 14).Value = -34401.0825   # N83: -35939 -> -34401.0825
# Row 122
$ws.Cells.Item(122, 8).Value = 31192   # H122: 20374.953 -> 31192
$ws.Cells.Item(122, 9).Value = 38444.223   # I122: 21169.176 -> 38444.223
$ws.Cells.Item(122, 10).Value = 14874.5   # J122: 16999.5 -> 14874.5
$ws.Cells.Item(122, 11).Value = 115332.669   # K122: 63507.528 -> 115332.669
$ws.Cells.Item(122, 12).Value = 44623.5   # L122: 50998.5 -> 44623.5
$ws.Cells.Item(122, 13).Value = -112882.669   # M122: -61057.528 -> -112882.669
$ws.Cells.Item(122, 14).Value = -49523.5   # N122: -55898.5 -> -49523.5
# Row 126
$ws.Cells.Item(126, 8).Value = 3100   # H126: 3050 -> 3100
$ws.Cells.Item(126, 10).Value = 0   # J126: 3000 -> 0
$ws.Cells.Item(126, 12).Value = 0   # L126: 9000 -> 0
$ws.Cells.Item(126, 14).Value = ""   # clear N126 (was -13940)

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 5844.0884   # H40: 6633.3 -> 5844.0884
$ws.Cells.Item(40, 9).Value = 3492.9333   # I40: 4363.091 -> 3492.9333
$ws.Cells.Item(40, 10).Value = 7700.263   # J40: 7947.6313 -> 7700.263
$ws.Cells.Item(40, 11).Value = 3492.9333   # K40: 4363.091 -> 3492.9333
$ws.Cells.Item(40, 12).Value = 7700.263   # L40: 7947.6313 -> 7700.263
$ws.Cells.Item(40, 13).Value = -3356.9333   # M40: -4227.091 -> -3356.9333
$ws.Cells.Item(40, 14).Value = -7972.263   # N40: -8219.631300000001 -> -7972.263
# Row 57
$ws.Cells.Item(57, 8).Value = 41046   # H57: 0 -> 41046
$ws.Cells.Item(57, 10).Value = 41046   # J57: 0 -> 41046
$ws.Cells.Item(57, 12).Value = 41046   # L57: 0 -> 41046
$ws.Cells.Item(57, 14).Value = -42178   # N57: None -> -42178
# Row 74
$ws.Cells.Item(74, 8).Value = 44125   # H74: 42800 -> 44125
$ws.Cells.Item(74, 9).Value = 38833.332   # I74: 38500 -> 38833.332
$ws.Cells.Item(74, 11).Value = 38833.332   # K74: 38500 -> 38833.332
$ws.Cells.Item(74, 13).Value = -37835.332   # M74: -37502 -> -37835.332
# Row 77
$ws.Cells.Item(77, 8).Value = 44125   # H77: 42800 -> 44125
$ws.Cells.Item(77, 9).Value = 38833.332   # I77: 38500 -> 38833.332
$ws.Cells.Item(77, 11).Value = 116499.996   # K77: 115500 -> 116499.996
$ws.Cells.Item(77, 13).Value = -111507.996   # M77: -110508 -> -111507.996
# Row 93
$ws.Cells.Item(93, 8).Value = 3885.5908   # H93: 3760.1304 -> 3885.5908
$ws.Cells.Item(93, 9).Value = 1639.6364   # I93: 1586.3334 -> 1639.6364
$ws.Cells.Item(93, 11).Value = 1639.6364   # K93: 1586.3334 -> 1639.6364
$ws.Cells.Item(93, 13).Value = -391.6364000000001   # M93: -338.3334 -> -391.6364000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 2236.6155   # H81: 2098.2856 -> 2236.6155
$ws.Cells.Item(81, 10).Value = 3987.8333   # J81: 3461 -> 3987.8333
$ws.Cells.Item(81, 12).Value = 7975.6666   # L81: 6922 -> 7975.6666
$ws.Cells.Item(81, 14).Value = -10097.6666   # N81: -9044 -> -10097.6666
# Row 84
$ws.Cells.Item(84, 8).Value = 2236.6155   # H84: 2098.2856 -> 2236.6155
$ws.Cells.Item(84, 10).Value = 3987.8333   # J84: 3461 -> 3987.8333
$ws.Cells.Item(84, 12).Value = 39878.333   # L84: 34610 -> 39878.333
$ws.Cells.Item(84, 14).Value = -50486.333   # N84: -45218 -> -50486.333
# Row 132
$ws.Cells.Item(132, 8).Value = 1591.32   # H132: 1624.2916 -> 1591.32
$ws.Cells.Item(132, 9).Value = 1588.1765   # I132: 1637.4375 -> 1588.1765
$ws.Cells.Item(132, 11).Value = 4764.529500000001   # K132: 4912.3125 -> 4764.529500000001
$ws.Cells.Item(132, 13).Value = -2234.529500000001   # M132: -2382.3125 -> -2234.529500000001
